$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of COVID overview data to append (2022-01-03 .. 2022-01-10)
$newRows = @(
    @{ Row = 510; Date = "2022-01-03"; CumCases = 13422815; NewCases = 157758; NewDeaths = 42;  CumDeaths = 148893 },
    @{ Row = 511; Date = "2022-01-04"; CumCases = 13641520; NewCases = 218724; NewDeaths = 48;  CumDeaths = 148941 },
    @{ Row = 512; Date = "2022-01-05"; CumCases = 13835334; NewCases = 194747; NewDeaths = 334; CumDeaths = 149284 },
    @{ Row = 513; Date = "2022-01-06"; CumCases = 14015065; NewCases = 179756; NewDeaths = 231; CumDeaths = 149515 },
    @{ Row = 514; Date = "2022-01-07"; CumCases = 14193228; NewCases = 178250; NewDeaths = 229; CumDeaths = 149744 },
    @{ Row = 515; Date = "2022-01-08"; CumCases = 14333794; NewCases = 146390; NewDeaths = 313; CumDeaths = 150057 },
    @{ Row = 516; Date = "2022-01-09"; CumCases = 14475192; NewCases = 141472; NewDeaths = 97;  CumDeaths = 150154 },
    @{ Row = 517; Date = "2022-01-10"; CumCases = 14617314; NewCases = 142224; NewDeaths = 77;  CumDeaths = 150230 }
)

# Column A holds date-like strings ("YYYY-MM-DD"). Force Text format on the
# whole block first so Excel doesn't auto-convert them into date serial
# numbers when assigned below.
$colA = $ws.Range("A510:A517")
$colA.NumberFormat = "@"

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = "overview"
    $ws.Cells.Item($row, 3).Value = "K02000001"
    $ws.Cells.Item($row, 4).Value = "United Kingdom"
    $ws.Cells.Item($row, 5).Value = $r.CumCases
    $ws.Cells.Item($row, 6).Value = $r.NewCases
    $ws.Cells.Item($row, 7).Value = $r.NewDeaths
    $ws.Cells.Item($row, 8).Value = $r.CumDeaths
}

# Restore the default "Normal" style on column A so the new rows carry no
# stray cell formatting (matching the rest of the sheet, which is unstyled).
$colA.Style = "Normal"
